$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Sprint 3"
$ws.Range("B4").Value = 41756
$ws.Range("C4").Value = 41756
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = 4

$ws.Range("B2:C2").Copy()
$ws.Range("B4:C4").PasteSpecial(-4122)

$ws.Range("E6").Select()
